$d = $word.ActiveDocument

# Revert "Get shortest path": remove the two paragraphs describing the
# txt-file example ("Bv. het zal een txt bestand inlezen, ..." and
# "Hoe langer de string zal zijn, ...") that were inserted before the
# "Het zal een MCMC structuur..." paragraph.
$startRange = $d.Content
$startRange.Find.Execute("Bv. het zal een txt bestand inlezen")
$startPos = $startRange.Start

$endRange = $d.Content
$endRange.Find.Execute("Het zal een MCMC structuur")
$endPos = $endRange.Start

$delRange = $d.Range($startPos, $endPos)
$delRange.Delete()

# Revert category change: "Goud" -> "Platina"
$d.Content.Find.Execute("Goud", $false, $false, $false, $false, $false, $true, 1, $false, "Platina", 2)
